$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) rows 2-498 all currently hold the serial date 45182
# (2023-09-13) and need to be updated to 45184 (2023-09-15).
$ws.Range("C2:C498").Value = 45184
